$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.848.00"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.58"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7216"
$ws.Range("E5").Value = "  -2.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.96"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3149"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07491"
$ws.Range("E9").Value = "  +3.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.59"
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08195"
$ws.Range("E11").Value = "  -2.52%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7451"
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.884.47"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.327"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.52"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.942.85"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.011"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.33"
$ws.Range("E18").Value = "  +1.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007922"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.49"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.154.05"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.723"
$ws.Range("E24").Value = "  -3.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.255"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1500"
$ws.Range("E26").Value = "  -3.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.56"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.59"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.011"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.435"
$ws.Range("E30").Value = "  -4.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.537"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.527"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.193"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05412"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7373"
$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.008"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.700"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01923"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4455"
$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8887"
$ws.Range("E42").Value = "  +3.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.013"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.60"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.040.44"
$ws.Range("E46").Value = "  -6.76%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.82"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.473"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.690"
$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.814"
$ws.Range("E50").Value = "  -1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.049.68"
$ws.Range("E51").Value = "  +1.33%  "
